# Auto-generated: applies cryptos.xlsx price/volume refresh + two coin swaps
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.594.65'
$ws.Range("E2").Value = '  -1.64%  '
$ws.Range("D3").Value = '1.667.81'
$ws.Range("E3").Value = '  -3.17%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''215.04'
$ws.Range("E5").Value = '  -1.73%  '
$ws.Range("E6").Value = '  -2.44%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '''23.69'
$ws.Range("E8").Value = '  -1.61%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("E10").Value = '  -1.54%  '
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("D12").Value = '1.902.56'
$ws.Range("E12").Value = '  -3.37%  '
$ws.Range("D13").Value = '1.666.24'
$ws.Range("E13").Value = '  -3.39%  '
$ws.Range("E14").Value = '  -3.08%  '
$ws.Range("D15").Value = '''0.562'
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").Value = '''66.24'
$ws.Range("E16").Value = '  -1.80%  '
$ws.Range("D17").Value = '27.604.21'
$ws.Range("E17").Value = '  -1.50%  '
$ws.Range("D18").Value = '''243.57'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("E19").Value = '  -3.46%  '
$ws.Range("D20").Value = '''7.61'
$ws.Range("E20").Value = '  -3.67%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("E22").Value = '  -3.10%  '
$ws.Range("E23").Value = '  -3.99%  '
$ws.Range("E24").Value = '  -4.34%  '
$ws.Range("D25").Value = '''146.72'
$ws.Range("E25").Value = '  -1.37%  '
$ws.Range("D26").Value = '''7.20'
$ws.Range("E26").Value = '  -3.99%  '
$ws.Range("E27").Value = '  -1.43%  '
$ws.Range("E30").Value = '  +2.77%  '
$ws.Range("D31").Value = '''0.0501'
$ws.Range("E31").Value = '  -1.81%  '
$ws.Range("E32").Value = '  -2.51%  '
$ws.Range("D33").Value = '1.465.97'
$ws.Range("E33").Value = '  -1.94%  '
$ws.Range("E34").Value = '  -4.67%  '
$ws.Range("D35").Value = '''1.57'
$ws.Range("E35").Value = '  -5.07%  '
$ws.Range("D36").Value = '''2.38'
$ws.Range("E36").Value = '  -1.68%  '
$ws.Range("D37").Value = '''0.931'
$ws.Range("E37").Value = '  -2.26%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '''0.575'
$ws.Range("E38").Value = '  -5.09%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.0172'
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D40").Value = '''69.48'
$ws.Range("E40").Value = '  -1.46%  '
$ws.Range("E41").Value = '  -5.18%  '
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.41'
$ws.Range("E43").Value = '  -7.00%  '
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").Value = '''2.22'
$ws.Range("E44").Value = '  -3.80%  '
$ws.Range("D45").Value = '1.810.74'
$ws.Range("E45").Value = '  -3.35%  '
$ws.Range("D46").Value = '''0.788'
$ws.Range("E46").Value = '  -0.96%  '
$ws.Range("D47").Value = '''1.72'
$ws.Range("E47").Value = '  -3.11%  '
$ws.Range("D48").Value = '''89.35'
$ws.Range("E48").Value = '  -1.69%  '
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("E50").Value = '  -2.14%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''7.88'
$ws.Range("E51").Value = '  -3.95%  '
